$wb = $excel.ActiveWorkbook

# --- Sheet "HLD" (sheet1) ---
$ws1 = $wb.Worksheets.Item("HLD")

# Row3: status Not Started -> Completed (new shared string "Completed")
$ws1.Range("C3").Value = "Completed"

# Row3: G3 comment (new shared string "Material UI")
$ws1.Range("G3").Value = "Material UI"

# --- Sheet "Landing Page" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Landing Page")

# New Row5 A5: (new shared string "Body")
$ws2.Range("A5").Value = "Body"

# Row4 (sheet1): fix typo in Redux description (new shared string, fixed text)
$ws1.Range("B4").Value = "Check which Redux library to use. Learn the differences and advantages."

# New Row5 B5 (sheet2): (new shared string "Add restaurants after the image.")
$ws2.Range("B5").Value = "Add restaurants after the image."

# Row2 D2: TODAY() formula -> static value (freeze current date value)
$ws1.Range("D2").ClearContents()
$ws1.Range("D2").Value = 44996

# Row3: D3 freeze date; E3 end date
$ws1.Range("D3").ClearContents()
$ws1.Range("D3").Value = 44996
$ws1.Range("E3").Value = 44997

# Row4: clear D4 (remove formula & value)
$ws1.Range("D4").ClearContents()

$ws1.Range("B4").Select()

# Row4 (sheet2): status Not Started -> In Progress; D4 gets a date value
$ws2.Range("C4").Value = "In Progress"
$ws2.Range("D4").Value = 44996

# New Row5 C5 (sheet2): In Progress
$ws2.Range("C5").Value = "In Progress"

$ws2.Range("B5").Select()
